$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number but must remain stored as
# text (matching the source data which keeps prices as literal strings).
# Force "Text" number format first so Excel does not auto-convert the string
# we assign into a numeric value, then restore the default "Normal" style so
# no visible formatting change is left behind.
$textCells = @('D4', 'D5', 'D6', 'D8', 'D9', 'D12', 'D13', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D30', 'D31', 'D32', 'D36', 'D40', 'D41', 'D42', 'D43', 'D46', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.987.61'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.914.95'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '590.03'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = '145.01'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.505'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').Value = '6.95'
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = '0.0000225'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '33.46'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '3.394.62'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '60.828.56'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '6.68'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '2.913.59'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').Value = '434.57'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '13.37'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '0.675'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = '7.10'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').Value = '81.50'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').Value = '10.94'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').Value = '11.78'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +4.66%  '
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').Value = '6.98'
$ws.Range('E30').Value = '  -3.83%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '26.49'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.109'
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').Value = '0.0₃0869'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').Value = '5.61'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('E39').Value = '  -2.79%  '
$ws.Range('D40').Value = '8.57'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('D41').Value = '41.99'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').Value = '0.288'
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('D43').Value = '377.01'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').Value = '2.687.49'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').Value = '132.80'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D48').Value = '23.84'
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('E51').Value = '  -0.90%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
